# Update workbook with corrected forecast output.
#
# Sheet1 ("Sheet1" -> "Sales vs PO") gains a new "Order Week" column (the
# original "ds" dates) while the "ds" column itself is rolled forward one
# week and the PO quantity column becomes a (zeroed) forward-looking
# forecast column. Three brand-new summary sheets are appended:
# "Weekly Growth", "Volume Insights" and "Prediction Info".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1: rename + restructure
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

# Remember the original "ds" dates (column A) before we shift anything
# around (column C / PO_Requested_Qty data just rides along with the
# column insert below, so it doesn't need to be cached).
$lastRow = 16
$origDs = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $origDs[$r] = $ws1.Cells.Item($r, 1).Value2
}

# Insert a new column at C - this pushes the old PO_Requested_Qty column
# (and its data) from C to D, inheriting formatting from the left.
$ws1.Columns.Item(3).Insert()

# Header for the newly inserted column.
$ws1.Range("C1").Value2 = "Order Week"
$ws1.Range("A1").Copy()
$ws1.Range("C1").PasteSpecial(-4122)
$ws1.Range("C1").Value2 = "Order Week"

# Populate the new "Order Week" column with the original "ds" dates, using
# the same date-formatted style as column A (copy format only, then write
# the value back so it doesn't inherit column A's post-shift date).
for ($r = 2; $r -le $lastRow; $r++) {
    $ws1.Range("A" + $r).Copy()
    $ws1.Range("C" + $r).PasteSpecial(-4122)
    $ws1.Range("C" + $r).Value2 = $origDs[$r]
}

# Roll the "ds" column forward one week (6 days) and zero out the
# PO_Requested_Qty forecast column (now column D).
for ($r = 2; $r -le $lastRow; $r++) {
    $ws1.Range("A" + $r).Value2 = $origDs[$r] + 6
    $ws1.Range("D" + $r).Value2 = 0
}

# ---------------------------------------------------------------------
# Sheet2: "Weekly Growth"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"

$ws2.Range("A1").Value2 = "ds"
$ws2.Range("B1").Value2 = "PO_Requested_Qty"
$ws2.Range("C1").Value2 = "Growth%"
$ws1.Range("A1:C1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)
$ws2.Range("A1").Value2 = "ds"
$ws2.Range("B1").Value2 = "PO_Requested_Qty"
$ws2.Range("C1").Value2 = "Growth%"

$growthDs = @(45558, 45565, 45593, 45628)
$growthPo = @(720, 24, 12, 60)
$growthPct = @(0, -96.66666666666667, -50, 400)

for ($i = 0; $i -lt $growthDs.Length; $i++) {
    $r = $i + 2
    $ws2.Range("A" + $r).Value2 = $growthDs[$i]
    $ws2.Range("B" + $r).Value2 = $growthPo[$i]
    $ws2.Range("C" + $r).Value2 = $growthPct[$i]
    $ws1.Range("A2").Copy()
    $ws2.Range("A" + $r).PasteSpecial(-4122)
    $ws2.Range("A" + $r).Value2 = $growthDs[$i]
}

# ---------------------------------------------------------------------
# Sheet3: "Volume Insights"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"

$ws3.Range("A1").Value2 = "Total_PO_Quantity"
$ws3.Range("B1").Value2 = "Average_PO_Quantity"
$ws3.Range("C1").Value2 = "Max_PO_Quantity"
$ws3.Range("D1").Value2 = "Min_PO_Quantity"
$ws1.Range("A1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)
$ws3.Range("A1").Value2 = "Total_PO_Quantity"
$ws3.Range("B1").Value2 = "Average_PO_Quantity"
$ws3.Range("C1").Value2 = "Max_PO_Quantity"
$ws3.Range("D1").Value2 = "Min_PO_Quantity"

$ws3.Range("A2").Value2 = 816
$ws3.Range("B2").Value2 = 204
$ws3.Range("C2").Value2 = 720
$ws3.Range("D2").Value2 = 12

# ---------------------------------------------------------------------
# Sheet4: "Prediction Info"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"

$ws4.Range("A1").Value2 = "Predicted_Next_Week_PO_Quantity"
$ws1.Range("A1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)
$ws4.Range("A1").Value2 = "Predicted_Next_Week_PO_Quantity"
$ws4.Range("A2").Value2 = 0

# ---------------------------------------------------------------------
# Page margins for the new sheets, matching the workbook-wide default
# (0.75in/0.75in/1in/1in/0.5in/0.5in == 54/54/72/72/36/36 pt).
# ---------------------------------------------------------------------
foreach ($sheet in @($ws2, $ws3, $ws4)) {
    $sheet.PageSetup.LeftMargin = 54
    $sheet.PageSetup.RightMargin = 54
    $sheet.PageSetup.TopMargin = 72
    $sheet.PageSetup.BottomMargin = 72
    $sheet.PageSetup.HeaderMargin = 36
    $sheet.PageSetup.FooterMargin = 36
}

# Leave the first sheet active/selected, matching the original workbook view.
$ws1.Select()
